$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders (permutes) the data rows 4-14 of the sheet (each row keeps all
# of its original column values together, rows are just redistributed to different
# row numbers). Column indices: A=1 .. T=20.
$firstCol = 1
$lastCol = 20

# Mapping: new row number -> original row number that its data comes from.
$rowMap = @{
    4  = 13
    5  = 12
    6  = 14
    7  = 10
    8  = 11
    9  = 4
    10 = 8
    11 = 7
    12 = 9
    13 = 5
    14 = 6
}

# Snapshot the original values of rows 4-14 (all columns) before overwriting anything,
# since several rows are sources for more than one destination and some rows are both
# a source and a destination.
$snapshot = @{}
foreach ($srcRow in 4..14) {
    $rowValues = @{}
    foreach ($col in $firstCol..$lastCol) {
        $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value()
    }
    $snapshot[$srcRow] = $rowValues
}

# Write the permuted data back into the sheet.
foreach ($destRow in 4..14) {
    $srcRow = $rowMap[$destRow]
    $rowValues = $snapshot[$srcRow]
    foreach ($col in $firstCol..$lastCol) {
        $ws.Cells.Item($destRow, $col).Value = $rowValues[$col]
    }
}
